$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 39

$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-01-23"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = "18:28:23"
$ws.Cells.Item($row, 3).Value = "Thursday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "03"
$ws.Cells.Item($row, 4).Style = "Normal"

$ws.Cells.Item($row, 5).Value = 126299
$ws.Cells.Item($row, 6).Value = 142140
$ws.Cells.Item($row, 7).Value = 168545
$ws.Cells.Item($row, 8).Value = 158651
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 142925
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 192285
$ws.Cells.Item($row, 14).Value = 115672
$ws.Cells.Item($row, 15).Value = 45592
$ws.Cells.Item($row, 16).Value = 28443
$ws.Cells.Item($row, 17).Value = 65521
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 47613
$ws.Cells.Item($row, 20).Value = -1
